{"js": "// Fix \"Convolution Neural Networks\" -> \"Convolutional Neural Networks\" in the\n// abstract paragraph (\"The rise of Deep learning has laid the path ...\").\n//\n// This mirrors what Word itself does when a user clicks right after\n// \"Convolution\" and types \"al\": the run is split and the \"_GoBack\" bookmark\n// (which always marks the location of the most recent edit) is relocated to\n// sit right after the newly inserted text.\n\nconst body = context.document.body;\n\n// Step 1: insert the missing \"al\" so \"Convolution\" becomes \"Convolutional\"\n// (only one occurrence of this phrase exists in the document).\nconst target = body.search(\"especially Convolution\", { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nconst convolutionRange = target.items[0];\nconvolutionRange.insertText(\"al\", Word.InsertLocation.end);\nawait context.sync();\n\n// Step 2: the \"_GoBack\" bookmark used to sit at the very end of the\n// paragraph; Word always moves it to the place of the latest edit, so drop\n// the old one.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Step 3: re-create \"_GoBack\" collapsed right after \"...especially\n// Convolutional\".\nconst target2 = body.search(\"especially Convolutional\", { matchCase: true });\ntarget2.load(\"text\");\nawait context.sync();\n\nconst convolutionalRange = target2.items[0];\nconst collapsedEnd = convolutionalRange.getRange(Word.RangeLocation.end);\ncollapsedEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Fix \"Convolution Neural Networks\" -> \"Convolutional Neural Networks\" in the\n# abstract paragraph (\"The rise of Deep learning has laid the path ...\").\n#\n# This mirrors what Word itself does when a user clicks right after\n# \"Convolution\" and types \"al\": the run is split and the \"_GoBack\" bookmark\n# (which always marks the location of the most recent edit) is relocated to\n# sit right after the newly inserted text.\n\n$d = $word.ActiveDocument\n\n# Step 1: insert the missing \"al\" so \"Convolution\" becomes \"Convolutional\"\n# (only one occurrence of this phrase exists in the document).\n$rngFind = $d.Content\n$rngFind.Find.MatchCase = $true\n$rngFind.Find.MatchWholeWord = $false\n$rngFind.Find.MatchWildcards = $false\n[void]$rngFind.Find.Execute(\"Convolution Neural N\", $false, $false, $false, $false, $false, $true, 1, $false, \"Convolutional Neural N\", 2)\n\n# Step 2: the \"_GoBack\" bookmark used to sit at the very end of the paragraph;\n# Word always moves it to the place of the latest edit, so drop the old one.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# Step 3: re-create \"_GoBack\" collapsed right after \"...especially Convolutional\".\n$rngMark = $d.Content\n$rngMark.Find.MatchCase = $true\n$rngMark.Find.MatchWholeWord = $false\n$rngMark.Find.MatchWildcards = $false\n[void]$rngMark.Find.Execute(\"especially Convolutional\")\n$rngMark.Collapse(0)\n[void]$d.Bookmarks.Add(\"_GoBack\", $rngMark)\n"}
